$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 (Sending cluster = "FAPs"/"sCs" originally) get re-derived values and
# rows 5-10 (new combinations with the "ECs" cluster) are added, per Dr Hou's
# advice on handling the natmi LR-pair edge-weight computation for 3 clusters
# (FAPs, sCs, ECs) instead of 2.

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Clcf1"
$ws.Cells.Item(2, 3).Value = "Cntfr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.9214586666666666
$ws.Cells.Item(2, 8).Value = 2.764376
$ws.Cells.Item(2, 9).Value = 0.08041853843186561
$ws.Cells.Item(2, 10).Value = 0.08041853843186561
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.04154133333333333
$ws.Cells.Item(2, 14).Value = 0.124624
$ws.Cells.Item(2, 15).Value = 0.005088925111573409
$ws.Cells.Item(2, 16).Value = 0.005088925111573409
$ws.Cells.Item(2, 17).Value = 0.03827862162488889
$ws.Cells.Item(2, 18).Value = 0.344507594624
$ws.Cells.Item(2, 19).Value = 0.0004092439196619522
$ws.Cells.Item(2, 20).Value = 0.0004092439196619522
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Clcf1"
$ws.Cells.Item(3, 3).Value = "Cntfr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.9214586666666666
$ws.Cells.Item(3, 8).Value = 2.764376
$ws.Cells.Item(3, 9).Value = 0.08041853843186561
$ws.Cells.Item(3, 10).Value = 0.08041853843186561
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.470089333333334
$ws.Cells.Item(3, 14).Value = 22.410268
$ws.Cells.Item(3, 15).Value = 0.915106043637582
$ws.Cells.Item(3, 16).Value = 0.9151060436375819
$ws.Cells.Item(3, 17).Value = 6.883378556974223
$ws.Cells.Item(3, 18).Value = 61.950407012768
$ws.Cells.Item(3, 19).Value = 0.07359149053950137
$ws.Cells.Item(3, 20).Value = 0.07359149053950137
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Clcf1"
$ws.Cells.Item(4, 3).Value = "Cntfr"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.9214586666666666
$ws.Cells.Item(4, 8).Value = 2.764376
$ws.Cells.Item(4, 9).Value = 0.08041853843186561
$ws.Cells.Item(4, 10).Value = 0.08041853843186561
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.6514553333333334
$ws.Cells.Item(4, 14).Value = 1.954366
$ws.Cells.Item(4, 15).Value = 0.07980503125084476
$ws.Cells.Item(4, 16).Value = 0.07980503125084475
$ws.Cells.Item(4, 17).Value = 0.6002891628462222
$ws.Cells.Item(4, 18).Value = 5.402602465616
$ws.Cells.Item(4, 19).Value = 0.006417803972702295
$ws.Cells.Item(4, 20).Value = 0.006417803972702294
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Clcf1"
$ws.Cells.Item(5, 3).Value = "Cntfr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.913147
$ws.Cells.Item(5, 8).Value = 5.739441
$ws.Cells.Item(5, 9).Value = 0.1669662363715809
$ws.Cells.Item(5, 10).Value = 0.1669662363715809
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.04154133333333333
$ws.Cells.Item(5, 14).Value = 0.124624
$ws.Cells.Item(5, 15).Value = 0.005088925111573409
$ws.Cells.Item(5, 16).Value = 0.005088925111573409
$ws.Cells.Item(5, 17).Value = 0.07947467724266667
$ws.Cells.Item(5, 18).Value = 0.715272095184
$ws.Cells.Item(5, 19).Value = 0.0008496786730562393
$ws.Cells.Item(5, 20).Value = 0.0008496786730562393
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Clcf1"
$ws.Cells.Item(6, 3).Value = "Cntfr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.913147
$ws.Cells.Item(6, 8).Value = 5.739441
$ws.Cells.Item(6, 9).Value = 0.1669662363715809
$ws.Cells.Item(6, 10).Value = 0.1669662363715809
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.470089333333334
$ws.Cells.Item(6, 14).Value = 22.410268
$ws.Cells.Item(6, 15).Value = 0.915106043637582
$ws.Cells.Item(6, 16).Value = 0.9151060436375819
$ws.Cells.Item(6, 17).Value = 14.29137899779867
$ws.Cells.Item(6, 18).Value = 128.622410980188
$ws.Cells.Item(6, 19).Value = 0.1527918119870547
$ws.Cells.Item(6, 20).Value = 0.1527918119870547
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Clcf1"
$ws.Cells.Item(7, 3).Value = "Cntfr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.913147
$ws.Cells.Item(7, 8).Value = 5.739441
$ws.Cells.Item(7, 9).Value = 0.1669662363715809
$ws.Cells.Item(7, 10).Value = 0.1669662363715809
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.6514553333333334
$ws.Cells.Item(7, 14).Value = 1.954366
$ws.Cells.Item(7, 15).Value = 0.07980503125084476
$ws.Cells.Item(7, 16).Value = 0.07980503125084475
$ws.Cells.Item(7, 17).Value = 1.246329816600667
$ws.Cells.Item(7, 18).Value = 11.216968349406
$ws.Cells.Item(7, 19).Value = 0.01332474571146994
$ws.Cells.Item(7, 20).Value = 0.01332474571146994
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Clcf1"
$ws.Cells.Item(8, 3).Value = "Cntfr"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.623680999999999
$ws.Cells.Item(8, 8).Value = 25.871043
$ws.Cells.Item(8, 9).Value = 0.7526152251965536
$ws.Cells.Item(8, 10).Value = 0.7526152251965536
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.04154133333333333
$ws.Cells.Item(8, 14).Value = 0.124624
$ws.Cells.Item(8, 15).Value = 0.005088925111573409
$ws.Cells.Item(8, 16).Value = 0.005088925111573409
$ws.Cells.Item(8, 17).Value = 0.3582392069813333
$ws.Cells.Item(8, 18).Value = 3.224152862832
$ws.Cells.Item(8, 19).Value = 0.003830002518855217
$ws.Cells.Item(8, 20).Value = 0.003830002518855217
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Clcf1"
$ws.Cells.Item(9, 3).Value = "Cntfr"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.623680999999999
$ws.Cells.Item(9, 8).Value = 25.871043
$ws.Cells.Item(9, 9).Value = 0.7526152251965536
$ws.Cells.Item(9, 10).Value = 0.7526152251965536
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 7.470089333333334
$ws.Cells.Item(9, 14).Value = 22.410268
$ws.Cells.Item(9, 15).Value = 0.915106043637582
$ws.Cells.Item(9, 16).Value = 0.9151060436375819
$ws.Cells.Item(9, 17).Value = 64.41966745216934
$ws.Cells.Item(9, 18).Value = 579.7770070695241
$ws.Cells.Item(9, 19).Value = 0.6887227411110259
$ws.Cells.Item(9, 20).Value = 0.6887227411110258
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Clcf1"
$ws.Cells.Item(10, 3).Value = "Cntfr"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 8.623680999999999
$ws.Cells.Item(10, 8).Value = 25.871043
$ws.Cells.Item(10, 9).Value = 0.7526152251965536
$ws.Cells.Item(10, 10).Value = 0.7526152251965536
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.6514553333333334
$ws.Cells.Item(10, 14).Value = 1.954366
$ws.Cells.Item(10, 15).Value = 0.07980503125084476
$ws.Cells.Item(10, 16).Value = 0.07980503125084475
$ws.Cells.Item(10, 17).Value = 5.617942980415333
$ws.Cells.Item(10, 18).Value = 50.561486823738
$ws.Cells.Item(10, 19).Value = 0.06006248156667252
$ws.Cells.Item(10, 20).Value = 0.06006248156667252
